$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header strings: volume/issue number and report week dates ---
$ws.Range("A8").Value = "Volume 32   Number  23"
$ws.Range("C9").Value = "Report Covering the Week  6/2/2025  Through  6/8/2025"

# --- Weekly crime-complaint figures (rows 15-33) ---
$ws.Range("C15").NumberFormat = "#,##0"
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = -50
$ws.Range("F15").Value = 6
$ws.Range("G15").Value = 8
$ws.Range("H15").Value = -25
$ws.Range("I15").Value = 34
$ws.Range("J15").Value = 29
$ws.Range("K15").Value = 17.241379310344
$ws.Range("L15").Value = 78.947368421052
$ws.Range("M15").Value = 54.545454545454
$ws.Range("N15").Value = -15
$ws.Range("C16").Value = 14
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 250
$ws.Range("F16").Value = 36
$ws.Range("G16").Value = 24
$ws.Range("H16").Value = 50
$ws.Range("I16").Value = 132
$ws.Range("J16").Value = 139
$ws.Range("K16").Value = -5.035971223021
$ws.Range("L16").Value = -4.347826086956
$ws.Range("M16").Value = -20.958083832335
$ws.Range("N16").Value = -77.358490566037
$ws.Range("D17").Value = 20
$ws.Range("E17").Value = 20
$ws.Range("F17").Value = 94
$ws.Range("G17").Value = 68
$ws.Range("H17").Value = 38.235294117647
$ws.Range("I17").Value = 439
$ws.Range("J17").Value = 427
$ws.Range("K17").Value = 2.810304449648
$ws.Range("L17").Value = 9.476309226932
$ws.Range("M17").Value = 97.747747747747
$ws.Range("N17").Value = -8.158995815899
$ws.Range("D18").Value = 8
$ws.Range("E18").Value = -62.5
$ws.Range("F18").Value = 19
$ws.Range("G18").Value = 35
$ws.Range("H18").Value = -45.714285714285
$ws.Range("I18").Value = 154
$ws.Range("J18").Value = 134
$ws.Range("K18").Value = 14.925373134328
$ws.Range("L18").Value = 1.986754966887
$ws.Range("M18").Value = -32.456140350877
$ws.Range("N18").Value = -89.516678012253
$ws.Range("C19").Value = 31
$ws.Range("D19").Value = 36
$ws.Range("E19").Value = -13.888888888888
$ws.Range("F19").Value = 110
$ws.Range("G19").Value = 126
$ws.Range("H19").Value = -12.698412698412
$ws.Range("I19").Value = 571
$ws.Range("J19").Value = 696
$ws.Range("K19").Value = -17.959770114942
$ws.Range("L19").Value = -16.763848396501
$ws.Range("M19").Value = 48.311688311688
$ws.Range("N19").Value = -14.903129657228
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 12
$ws.Range("E20").Value = -58.333333333333
$ws.Range("F20").Value = 14
$ws.Range("G20").Value = 41
$ws.Range("H20").Value = -65.853658536585
$ws.Range("I20").Value = 80
$ws.Range("J20").Value = 129
$ws.Range("K20").Value = -37.984496124031
$ws.Range("L20").Value = -55.05617977528
$ws.Range("M20").Value = -36.507936507936
$ws.Range("N20").Value = -96.217494089834
$ws.Range("C21").Value = 78
$ws.Range("D21").Value = 82
$ws.Range("E21").Value = -4.878048780487
$ws.Range("F21").Value = 279
$ws.Range("G21").Value = 302
$ws.Range("H21").Value = -7.615894039735
$ws.Range("I21").Value = 1412
$ws.Range("J21").Value = 1555
$ws.Range("K21").Value = -9.196141479099
$ws.Range("L21").Value = -10.802274162981
$ws.Range("M21").Value = 22.463139635732
$ws.Range("N21").Value = -73.691075088503
$ws.Range("C23").NumberFormat = "#,##0"
$ws.Range("C23").Value = 5
$ws.Range("D23").Value = 4
$ws.Range("E23").Value = 25
$ws.Range("F23").Value = 13
$ws.Range("G23").Value = 8
$ws.Range("H23").Value = 62.5
$ws.Range("I23").Value = 66
$ws.Range("J23").Value = 49
$ws.Range("K23").Value = 34.69387755102
$ws.Range("L23").Value = -7.042253521126
$ws.Range("M23").Value = 164
$ws.Range("C24").Value = 84
$ws.Range("D24").Value = 69
$ws.Range("E24").Value = 21.739130434782
$ws.Range("F24").Value = 330
$ws.Range("G24").Value = 301
$ws.Range("H24").Value = 9.634551495016
$ws.Range("I24").Value = 1880
$ws.Range("J24").Value = 1865
$ws.Range("K24").Value = 0.804289544235
$ws.Range("L24").Value = 4.67706013363
$ws.Range("M24").Value = 23.765635286372
$ws.Range("C25").Value = 55
$ws.Range("D25").Value = 36
$ws.Range("E25").Value = 52.777777777777
$ws.Range("F25").Value = 210
$ws.Range("H25").Value = 32.911392405063
$ws.Range("I25").Value = 1130
$ws.Range("J25").Value = 1029
$ws.Range("K25").Value = 9.815354713313
$ws.Range("L25").Value = 30.034522439585
$ws.Range("C26").Value = 35
$ws.Range("D26").Value = 40
$ws.Range("E26").Value = -12.5
$ws.Range("F26").Value = 139
$ws.Range("G26").Value = 152
$ws.Range("H26").Value = -8.552631578947
$ws.Range("I26").Value = 819
$ws.Range("J26").Value = 796
$ws.Range("K26").Value = 2.88944723618
$ws.Range("L26").Value = 7.339449541284
$ws.Range("M26").Value = -5.970149253731
$ws.Range("C27").NumberFormat = "#,##0"
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 3
$ws.Range("E27").Value = -66.666666666666
$ws.Range("F27").Value = 6
$ws.Range("G27").Value = 14
$ws.Range("H27").Value = -57.142857142857
$ws.Range("I27").Value = 42
$ws.Range("J27").Value = 52
$ws.Range("K27").Value = -19.230769230769
$ws.Range("L27").Value = 44.827586206896
$ws.Range("C28").Value = 2
$ws.Range("D28").Value = 5
$ws.Range("E28").Value = -60
$ws.Range("F28").Value = 12
$ws.Range("G28").Value = 17
$ws.Range("H28").Value = -29.411764705882
$ws.Range("I28").Value = 97
$ws.Range("J28").Value = 82
$ws.Range("K28").Value = 18.292682926829
$ws.Range("L28").Value = 4.301075268817
$ws.Range("L29").Value = -75
$ws.Range("N29").Value = -94
$ws.Range("L30").Value = -75
$ws.Range("N30").Value = -92.5
$ws.Range("G31").Value = 3
$ws.Range("J31").Value = 11
$ws.Range("K31").Value = -63.636363636363
$ws.Range("L31").Value = -33.333333333333
$ws.Range("D33").NumberFormat = "#,##0"
$ws.Range("D33").Value = 1
$ws.Range("E33").NumberFormat = "#,##0.0;`"-`"#,##0.0"
$ws.Range("E33").Value = -100
$ws.Range("G33").NumberFormat = "#,##0"
$ws.Range("G33").Value = 1
$ws.Range("H33").NumberFormat = "#,##0.0;`"-`"#,##0.0"
$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 4
$ws.Range("K33").Value = 0
